# Adds Problem2Bonus sheet with benchmarking statistics for problem 2,
# and tweaks the selection/active sheet state to match.

$wb = $excel.ActiveWorkbook

# --- Update the cursor/selection on the existing "Problem2" sheet ---
$ws2 = $wb.Worksheets.Item("Problem2")
$ws2.Range("D9").Select() | Out-Null

# --- Add the new "Problem2Bonus" worksheet after the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Problem2Bonus"

$colA = @(
    "time=150.98047800 seconds",
    "Memory Reads[2] = 121870480",
    "Memory Writes[2] = 11695138",
    "Basic Blocks[2] = 7594768",
    "Total Instructions[2] = 197563103",
    "Memory Reads[1] = 121545494",
    "Memory Writes[1] = 11695306",
    "Basic Blocks[1] = 6944212",
    "Total Instructions[1] = 195611910",
    "Memory Reads[3] = 121854939",
    "Memory Writes[3] = 11694888",
    "Basic Blocks[3] = 7564592",
    "Total Instructions[3] = 197471670",
    "Memory Reads[0] = 125361231",
    "Memory Writes[0] = 13353706",
    "Basic Blocks[0] = 9019806",
    "Total Instructions[0] = 209355661",
    "Total number of threads = 4"
)

$colF = @(
    "time=0.16943400 seconds",
    "Memory Reads[3] = 68104",
    "Memory Writes[3] = 1528",
    "Basic Blocks[3] = 133237",
    "Total Instructions[3] = 403338",
    "Memory Reads[2] = 123182",
    "Memory Writes[2] = 1521",
    "Basic Blocks[2] = 243396",
    "Total Instructions[2] = 733805",
    "Memory Reads[1] = 128701",
    "Memory Writes[1] = 1804",
    "Basic Blocks[1] = 253410",
    "Total Instructions[1] = 764766",
    "Memory Reads[0] = 349752",
    "Memory Writes[0] = 119591",
    "Basic Blocks[0] = 252405",
    "Total Instructions[0] = 1449429",
    "Total number of threads = 4"
)

$colK = @(
    "time=143.96135900 seconds",
    "Memory Reads[2] = 121907204",
    "Memory Writes[2] = 11695170",
    "Basic Blocks[2] = 7668105",
    "Total Instructions[2] = 197783103",
    "Memory Reads[3] = 121851383",
    "Memory Writes[3] = 11695014",
    "Basic Blocks[3] = 7557022",
    "Total Instructions[3] = 197449416",
    "Memory Reads[1] = 121595688",
    "Memory Writes[1] = 11695305",
    "Basic Blocks[1] = 7044600",
    "Total Instructions[1] = 195913074",
    "Memory Reads[0] = 125361271",
    "Memory Writes[0] = 13353728",
    "Basic Blocks[0] = 9019230",
    "Total Instructions[0] = 209354000",
    "Total number of threads = 4"
)

$colO = @(
    "time=0.17670500 seconds",
    "Memory Reads[1] = 136010",
    "Memory Writes[1] = 1808",
    "Basic Blocks[1] = 268022",
    "Total Instructions[1] = 808608",
    "Memory Reads[2] = 122385",
    "Memory Writes[2] = 1648",
    "Basic Blocks[2] = 241347",
    "Total Instructions[2] = 728118",
    "Memory Reads[3] = 82365",
    "Memory Writes[3] = 1639",
    "Basic Blocks[3] = 161317",
    "Total Instructions[3] = 488016",
    "Memory Reads[0] = 348213",
    "Memory Writes[0] = 119347",
    "Basic Blocks[0] = 249596",
    "Total Instructions[0] = 1440102",
    "Total number of threads = 4"
)

# Fill column-by-column (not row-by-row) so that shared-string entries are
# interned in the same order the original workbook used.
for ($i = 0; $i -lt $colA.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt $colF.Length; $i++) {
    $ws.Cells.Item($i + 1, 6).Value = $colF[$i]
}
for ($i = 0; $i -lt $colK.Length; $i++) {
    $ws.Cells.Item($i + 1, 11).Value = $colK[$i]
}
for ($i = 0; $i -lt $colO.Length; $i++) {
    $ws.Cells.Item($i + 1, 15).Value = $colO[$i]
}

$ws.Range("I20").Select() | Out-Null
